$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new year header 2022 in S4, formatted like R4 ---
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# --- Row 5 ---
$ws.Range("Q5").Value = 117.60684979252385
$ws.Range("R5").Value = 113.34848864817617
$ws.Range("A6").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 115.8

# --- Row 6 ---
$ws.Range("Q6").Value = 114.77319768114526
$ws.Range("R6").Value = 115.06069350712495
$ws.Range("A6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 115.2

# --- Row 7 ---
$ws.Range("Q7").Value = 116.40044011407315
$ws.Range("R7").Value = 114.29658549692938
$ws.Range("A6").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 115.4

# --- Row 8 ---
$ws.Range("Q8").Value = 117.53828537152096
$ws.Range("R8").Value = 113.75761785228545
$ws.Range("A6").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = 111.8

# --- Row 9 ---
$ws.Range("Q9").Value = 117.42206669681742
$ws.Range("R9").Value = 113.98264089946031
$ws.Range("A6").Copy()
$ws.Range("S9").PasteSpecial(-4122)
$ws.Range("S9").Value = 116.8

# --- Row 10 ---
$ws.Range("Q10").Value = 113.98326995089161
$ws.Range("R10").Value = 113.92720567782911
$ws.Range("A6").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("S10").Value = 108.2

# --- Row 11 ---
$ws.Range("Q11").Value = 123.488978736909
$ws.Range("R11").Value = 114.17226706705155
$ws.Range("A6").Copy()
$ws.Range("S11").PasteSpecial(-4122)
$ws.Range("S11").Value = 111

# --- Row 12 ---
$ws.Range("Q12").Value = 118.12340252754679
$ws.Range("R12").Value = 114.45153946490467
$ws.Range("A6").Copy()
$ws.Range("S12").PasteSpecial(-4122)
$ws.Range("S12").Value = 115.8

# --- Row 13 ---
$ws.Range("Q13").Value = 118.87059844457349
$ws.Range("R13").Value = 112.69493421065988
$ws.Range("A6").Copy()
$ws.Range("S13").PasteSpecial(-4122)
$ws.Range("S13").Value = 117.9

# --- Row 14 ---
$ws.Range("Q14").Value = 114.06377070452145
$ws.Range("R14").Value = 113.95067699644588
$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)
$ws.Range("S14").Value = 112.4

$excel.CutCopyMode = 0

# --- Selection update ---
$ws.Range("T4").Select()
